{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Insert the new \"Buses cannot have 0 capacity.\" restriction bullet right\n//        after \"Bus cannot start another trip before finishing the previous one.\",\n//        matching the style/numbering of its sibling bullets, and move the single\n//        \"_GoBack\" bookmark (Word's \"last edit\" marker) onto the end of this new\n//        sentence - exactly where Word itself leaves it after the newest edit.\nlet anchor = null;\nfor (const para of paragraphs.items) {\n  if (para.text === \"Bus cannot start another trip before finishing the previous one.\") {\n    anchor = para;\n    break;\n  }\n}\n\n// Word keeps only a single \"_GoBack\" bookmark in the whole document, so drop the\n// old one before planting the new one.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Type a sentinel character after the real sentence so the bookmark's collapsed\n// insertion point is never the very last offset of the paragraph (a boundary\n// position that is ambiguous with the start of the following paragraph); then\n// drop the bookmark right before the sentinel and remove the sentinel.\nconst newPara = anchor.insertParagraph(\"Buses cannot have 0 capacity.#\", Word.InsertLocation.after);\nawait context.sync();\n\nconst newParaContent = newPara.getRange(Word.RangeLocation.content);\nconst sentinelMatches = newParaContent.search(\"#\", { matchWildcards: false });\nsentinelMatches.load(\"items\");\nawait context.sync();\n\nconst sentinelRange = sentinelMatches.items[0];\nconst bookmarkRange = sentinelRange.getRange(Word.RangeLocation.start);\nbookmarkRange.insertBookmark(\"_GoBack\");\nsentinelRange.delete();\nawait context.sync();\n\n// --- 2. Normalise the \"Can book trips.\" paragraph so its text lives in a single\n//        run instead of being split across \"C\" / \"an book trips.\" runs.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet canBookTrips = null;\nfor (const para of paragraphs.items) {\n  if (para.text === \"Can book trips.\") {\n    canBookTrips = para;\n    break;\n  }\n}\ncanBookTrips.insertText(\"Can book trips.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Insert new restriction bullet \"Buses cannot have 0 capacity.\" right after\n#        \"Bus cannot start another trip before finishing the previous one.\" and leave\n#        the _GoBack bookmark collapsed at the end of the new sentence (Word moves the\n#        single \"last edit\" bookmark there automatically when new text is typed).\n$anchor = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $cand = $d.Paragraphs.Item($i)\n    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq \"Bus cannot start another trip before finishing the previous one.\") {\n        $anchor = $cand\n        break\n    }\n}\n\n$anchor.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item($anchor.Index + 1)\n# Type with a temporary trailing sentinel character so the true end-of-text position\n# isn't the (fragile) very-last offset of the paragraph while we drop the bookmark.\n$newPara.Range.Text = \"Buses cannot have 0 capacity.#\"\n$endPos = $newPara.Range.End - 2\n$bmRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n# Remove the temporary sentinel character now that the bookmark is anchored.\n$sentinel = $d.Range($endPos, $endPos + 1)\n$sentinel.Text = \"\"\n\n# --- 2. Normalise the \"Can book trips.\" paragraph so its text lives in a single run\n#        instead of being split across \"C\" / \"an book trips.\".\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $cand = $d.Paragraphs.Item($i)\n    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq \"Can book trips.\") {\n        $target = $cand\n        break\n    }\n}\n$r = $target.Range.Duplicate()\n[void]$r.MoveEnd(1, -1)\n$r.Delete()\n$r.InsertAfter(\"Can book trips.\")\n\nWrite-Output \"done\"\n"}
